$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("ZZ1")

$ws.Range('E2').Value = '2026-02-19 22:18:46'
$helper.Formula = '="70%"'
$helper.Copy()
$ws.Range('H2').PasteSpecial(-4163)
$ws.Range('I2').Value = '4.0 mm'
$ws.Range('E3').Value = '2026-02-19 22:18:49'
$ws.Range('I3').Value = '6.0 mm'
$ws.Range('E4').Value = '2026-02-19 22:18:52'
$ws.Range('J4').Value = '1010.2 hPa'
$ws.Range('E5').Value = '2026-02-19 22:18:54'
$ws.Range('I5').Value = '8.1 mm'
$ws.Range('O5').Value = '-6.1 °C'
$ws.Range('E6').Value = '2026-02-19 22:18:57'
$ws.Range('J6').Value = '1010.4 hPa'
$ws.Range('E7').Value = '2026-02-19 22:19:00'
$ws.Range('J7').Value = '1011.4 hPa'
$ws.Range('E8').Value = '2026-02-19 22:19:02'
$ws.Range('J8').Value = '1011.1 hPa'
$ws.Range('O8').Value = '9.8 °C'
$ws.Range('E9').Value = '2026-02-19 22:19:05'
$ws.Range('E10').Value = '2026-02-19 22:19:08'
$ws.Range('N10').Value = '3.4 °C 21:49 TU'
$ws.Range('O10').Value = '9.9 °C'
$ws.Range('E11').Value = '2026-02-19 22:19:10'
$ws.Range('O11').Value = '5.9 °C'
$ws.Range('E12').Value = '2026-02-19 22:19:13'
$helper.Formula = '="79%"'
$helper.Copy()
$ws.Range('H12').PasteSpecial(-4163)
$ws.Range('E13').Value = '2026-02-19 22:19:15'
$ws.Range('J13').Value = '1011.7 hPa'
$ws.Range('O13').Value = '4.5 °C'
$ws.Range('E14').Value = '2026-02-19 22:19:18'
$ws.Range('O14').Value = '13.1 °C'
$ws.Range('E15').Value = '2026-02-19 22:19:21'
$ws.Range('O15').Value = '9.7 °C'
$ws.Range('E16').Value = '2026-02-19 22:19:23'
$ws.Range('I16').Value = '11.4 mm'
$ws.Range('E17').Value = '2026-02-19 22:19:26'
$helper.Formula = '="79%"'
$helper.Copy()
$ws.Range('H17').PasteSpecial(-4163)
$ws.Range('E18').Value = '2026-02-19 22:19:29'
$helper.Formula = '="61%"'
$helper.Copy()
$ws.Range('H18').PasteSpecial(-4163)
$ws.Range('J18').Value = '1010.5 hPa'
$ws.Range('N18').Value = '5.4 °C 21:59 TU'
$ws.Range('O18').Value = '11.4 °C'
$ws.Range('E19').Value = '2026-02-19 22:19:31'
$ws.Range('E20').Value = '2026-02-19 22:19:34'
$ws.Range('E21').Value = '2026-02-19 22:19:37'
$ws.Range('J21').Value = '1011.8 hPa'
$ws.Range('K21').Value = '13.0 MJ/m2'
$ws.Range('E22').Value = '2026-02-19 22:19:40'
$ws.Range('L22').Value = '119.2 km/h - 336º 21:53 TU'
$ws.Range('E23').Value = '2026-02-19 22:19:43'
$ws.Range('G23').Value = '216 cm'
$ws.Range('I23').Value = '11.5 mm'
$ws.Range('E24').Value = '2026-02-19 22:19:45'
$ws.Range('J24').Value = '1015.3 hPa'
$ws.Range('E25').Value = '2026-02-19 22:19:48'
$ws.Range('I25').Value = '7.4 mm'
$ws.Range('E26').Value = '2026-02-19 22:19:51'
$ws.Range('J26').Value = '1010.2 hPa'
$ws.Range('L26').Value = '77.4 km/h - 315º 21:50 TU'
$ws.Range('E27').Value = '2026-02-19 22:19:53'
$helper.Formula = '="67%"'
$helper.Copy()
$ws.Range('H27').PasteSpecial(-4163)
$ws.Range('E28').Value = '2026-02-19 22:19:56'
$ws.Range('J28').Value = '1010.2 hPa'
$ws.Range('O28').Value = '9.1 °C'
$ws.Range('E29').Value = '2026-02-19 22:19:59'
$ws.Range('N29').Value = '4.5 °C 21:51 TU'
$ws.Range('O29').Value = '10.2 °C'
$ws.Range('E30').Value = '2026-02-19 22:20:01'
$ws.Range('J30').Value = '1010.4 hPa'
$ws.Range('O30').Value = '9.9 °C'
$ws.Range('E31').Value = '2026-02-19 22:20:04'
$ws.Range('J31').Value = '1009.8 hPa'
$ws.Range('E32').Value = '2026-02-19 22:20:06'
$ws.Range('E33').Value = '2026-02-19 22:20:09'
$ws.Range('J33').Value = '1011.2 hPa'
$ws.Range('E34').Value = '2026-02-19 22:20:12'
$ws.Range('E35').Value = '2026-02-19 22:20:15'
$ws.Range('J35').Value = '1016.8 hPa'
$ws.Range('L35').Value = '82.1 km/h - 269º 21:51 TU'
$ws.Range('O35').Value = '4.0 °C'
$ws.Range('E36').Value = '2026-02-19 22:20:18'
$ws.Range('J36').Value = '1010.6 hPa'
$ws.Range('O36').Value = '11.9 °C'
$ws.Range('E37').Value = '2026-02-19 22:20:20'
$ws.Range('J37').Value = '1011.7 hPa'
$ws.Range('E38').Value = '2026-02-19 22:20:23'
$ws.Range('E39').Value = '2026-02-19 22:20:26'
$helper.Formula = '="74%"'
$helper.Copy()
$ws.Range('H39').PasteSpecial(-4163)
$ws.Range('I39').Value = '5.0 mm'
$ws.Range('E40').Value = '2026-02-19 22:20:29'
$ws.Range('J40').Value = '1013.0 hPa'
$ws.Range('E41').Value = '2026-02-19 22:20:31'
$helper.Formula = '="40%"'
$helper.Copy()
$ws.Range('H41').PasteSpecial(-4163)
$ws.Range('J41').Value = '1013.3 hPa'
$ws.Range('E42').Value = '2026-02-19 22:20:34'
$ws.Range('N42').Value = '5.9 °C 21:57 TU'
$ws.Range('O42').Value = '10.9 °C'
$ws.Range('E43').Value = '2026-02-19 22:20:37'
$ws.Range('E44').Value = '2026-02-19 22:20:39'
$ws.Range('I44').Value = '10.1 mm'
$ws.Range('E45').Value = '2026-02-19 22:20:42'
$helper.Formula = '="83%"'
$helper.Copy()
$ws.Range('H45').PasteSpecial(-4163)
$ws.Range('I45').Value = '3.5 mm'
$ws.Range('J45').Value = '1016.0 hPa'
$ws.Range('O45').Value = '2.7 °C'
$ws.Range('E46').Value = '2026-02-19 22:20:45'
$ws.Range('J46').Value = '1016.2 hPa'
$helper.ClearContents()
